$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "arrival time (local time)" column (M),
# shifting M:V to N:W.
$ws.Columns("M").Insert()

# New header for the inserted column.
$ws.Range("M1").Value = "fuel_remaining(liters)"

# Fuel remaining (liters) values for each leg row.
$ws.Range("M2").Value = 9
$ws.Range("M4").Value = 202
$ws.Range("M5").Value = 65
$ws.Range("M8").Value = 1118
$ws.Range("M9").Value = 836
$ws.Range("M11").Value = 716
$ws.Range("M12").Value = 342
$ws.Range("M15").Value = 727
$ws.Range("M16").Value = 342
$ws.Range("M18").Value = 9

# Update the "description of work" text for the Summit refuel leg (now column W).
$ws.Range("W4").Value = "refuel at Summit 2 x 200 liters. basic maintanance at Summit"
